$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01818866666666667
$ws.Range("H2").Value = 0.054566
$ws.Range("I2").Value = 0.006403810693375696
$ws.Range("J2").Value = 0.006403810693375696
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.175215
$ws.Range("N2").Value = 6.525645
$ws.Range("O2").Value = 0.1204575739226287
$ws.Range("P2").Value = 0.1204575739226287
$ws.Range("Q2").Value = 0.03956426056333334
$ws.Range("R2").Value = 0.35607834507
$ws.Range("S2").Value = 0.0007713874999838233
$ws.Range("T2").Value = 0.0007713874999838232

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01818866666666667
$ws.Range("H3").Value = 0.054566
$ws.Range("I3").Value = 0.006403810693375696
$ws.Range("J3").Value = 0.006403810693375696
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.871597
$ws.Range("N3").Value = 32.614791
$ws.Range("O3").Value = 0.6020398899807737
$ws.Range("P3").Value = 0.6020398899807737
$ws.Range("Q3").Value = 0.1977398539673333
$ws.Range("R3").Value = 1.779658685706
$ws.Range("S3").Value = 0.003855349485297606
$ws.Range("T3").Value = 0.003855349485297606

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr4"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01818866666666667
$ws.Range("H4").Value = 0.054566
$ws.Range("I4").Value = 0.006403810693375696
$ws.Range("J4").Value = 0.006403810693375696
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.428447
$ws.Range("N4").Value = 1.285341
$ws.Range("O4").Value = 0.0237262459915128
$ws.Range("P4").Value = 0.0237262459915128
$ws.Range("Q4").Value = 0.007792879667333334
$ws.Range("R4").Value = 0.070135917006
$ws.Range("S4").Value = 0.0001519383877941119
$ws.Range("T4").Value = 0.0001519383877941119

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01818866666666667
$ws.Range("H5").Value = 0.054566
$ws.Range("I5").Value = 0.006403810693375696
$ws.Range("J5").Value = 0.006403810693375696
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.202622333333333
$ws.Range("N5").Value = 9.607866999999999
$ws.Range("O5").Value = 0.1773526370789838
$ws.Range("P5").Value = 0.1773526370789838
$ws.Range("Q5").Value = 0.05825143008022222
$ws.Range("R5").Value = 0.524262870722
$ws.Range("S5").Value = 0.001135732713824775
$ws.Range("T5").Value = 0.001135732713824775

$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Lgr4"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01818866666666667
$ws.Range("H6").Value = 0.054566
$ws.Range("I6").Value = 0.006403810693375696
$ws.Range("J6").Value = 0.006403810693375696
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.223995666666667
$ws.Range("N6").Value = 3.671987
$ws.Range("O6").Value = 0.06778159791031105
$ws.Range("P6").Value = 0.06778159791031105
$ws.Range("Q6").Value = 0.02226284918244445
$ws.Range("R6").Value = 0.200365642642
$ws.Range("S6").Value = 0.0004340605215121416
$ws.Range("T6").Value = 0.0004340605215121416

$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Lgr4"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01818866666666667
$ws.Range("H7").Value = 0.054566
$ws.Range("I7").Value = 0.006403810693375696
$ws.Range("J7").Value = 0.006403810693375696
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1560576666666667
$ws.Range("N7").Value = 0.468173
$ws.Range("O7").Value = 0.008642055115789912
$ws.Range("P7").Value = 0.008642055115789912
$ws.Range("Q7").Value = 0.002838480879777778
$ws.Range("R7").Value = 0.025546327918
$ws.Range("S7").Value = 0.00005534208496323757
$ws.Range("T7").Value = 0.00005534208496323757

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rspo3"
$ws.Range("C8").Value = "Lgr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.822099333333334
$ws.Range("H8").Value = 8.466298
$ws.Range("I8").Value = 0.9935961893066243
$ws.Range("J8").Value = 0.9935961893066244
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.175215
$ws.Range("N8").Value = 6.525645
$ws.Range("O8").Value = 0.1204575739226287
$ws.Range("P8").Value = 0.1204575739226287
$ws.Range("Q8").Value = 6.138672801356668
$ws.Range("R8").Value = 55.24805521221
$ws.Range("S8").Value = 0.1196861864226449
$ws.Range("T8").Value = 0.1196861864226449

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rspo3"
$ws.Range("C9").Value = "Lgr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.822099333333334
$ws.Range("H9").Value = 8.466298
$ws.Range("I9").Value = 0.9935961893066243
$ws.Range("J9").Value = 0.9935961893066244
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.871597
$ws.Range("N9").Value = 32.614791
$ws.Range("O9").Value = 0.6020398899807737
$ws.Range("P9").Value = 0.6020398899807737
$ws.Range("Q9").Value = 30.68072664596867
$ws.Range("R9").Value = 276.126539813718
$ws.Range("S9").Value = 0.5981845404954761
$ws.Range("T9").Value = 0.5981845404954762

$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Rspo3"
$ws.Range("C10").Value = "Lgr4"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.822099333333334
$ws.Range("H10").Value = 8.466298
$ws.Range("I10").Value = 0.9935961893066243
$ws.Range("J10").Value = 0.9935961893066244
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.428447
$ws.Range("N10").Value = 1.285341
$ws.Range("O10").Value = 0.0237262459915128
$ws.Range("P10").Value = 0.0237262459915128
$ws.Range("Q10").Value = 1.209119993068667
$ws.Range("R10").Value = 10.882079937618
$ws.Range("S10").Value = 0.02357430760371869
$ws.Range("T10").Value = 0.02357430760371869

$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Rspo3"
$ws.Range("C11").Value = "Lgr4"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.822099333333334
$ws.Range("H11").Value = 8.466298
$ws.Range("I11").Value = 0.9935961893066243
$ws.Range("J11").Value = 0.9935961893066244
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.202622333333333
$ws.Range("N11").Value = 9.607866999999999
$ws.Range("O11").Value = 0.1773526370789838
$ws.Range("P11").Value = 0.1773526370789838
$ws.Range("Q11").Value = 9.038118351818444
$ws.Range("R11").Value = 81.34306516636599
$ws.Range("S11").Value = 0.176216904365159
$ws.Range("T11").Value = 0.176216904365159

$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Rspo3"
$ws.Range("C12").Value = "Lgr4"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.822099333333334
$ws.Range("H12").Value = 8.466298
$ws.Range("I12").Value = 0.9935961893066243
$ws.Range("J12").Value = 0.9935961893066244
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.223995666666667
$ws.Range("N12").Value = 3.671987
$ws.Range("O12").Value = 0.06778159791031105
$ws.Range("P12").Value = 0.06778159791031105
$ws.Range("Q12").Value = 3.454237354902889
$ws.Range("R12").Value = 31.088136194126
$ws.Range("S12").Value = 0.06734753738879891
$ws.Range("T12").Value = 0.06734753738879891

$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Rspo3"
$ws.Range("C13").Value = "Lgr4"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.822099333333334
$ws.Range("H13").Value = 8.466298
$ws.Range("I13").Value = 0.9935961893066243
$ws.Range("J13").Value = 0.9935961893066244
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1560576666666667
$ws.Range("N13").Value = 0.468173
$ws.Range("O13").Value = 0.008642055115789912
$ws.Range("P13").Value = 0.008642055115789912
$ws.Range("Q13").Value = 0.4404102370615555
$ws.Range("R13").Value = 3.963692133553999
$ws.Range("S13").Value = 0.008586713030826674
$ws.Range("T13").Value = 0.008586713030826676

Write-Output "Applied updated TPM data"